$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 2030 electricity demand values (Value column, E21:E23)
$ws.Range("E21").Value = 500
$ws.Range("E22").Value = 500
$ws.Range("E23").Value = 500

# Remove the existing AutoFilter so it can be re-applied over the new, larger range
$ws.AutoFilterMode = $false

# Re-apply the AutoFilter over the full data range (now A1:E31) filtering Year (col 4) to 2030
$ws.Range("A1:E31").AutoFilter(4, @("2030"), 7)

# Keep the hidden _xlnm._FilterDatabase defined name in sync with the new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "demand!_FilterDatabase") {
        $n.RefersTo = "=demand!`$A`$1:`$E`$31"
    }
}

# Update the active cell / selection to match the author's last position
$ws.Range("E37").Select()
